$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference a plain-styled, never-edited cell to normalize style index
# after using the quote-prefix trick to force text for numeric-looking values.
$plainStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = '68.656.42'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '2.456.79'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'557.62"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = "'161.57"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.510"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("D12").Value = "'4.84"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = '68.553.44'
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").Value = "'23.39"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").Value = "'10.57"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = '  -3.61%  '
$ws.Range("D17").Value = "'334.78"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").Value = "'3.76"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = "'66.28"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").Value = "'3.63"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("D24").Value = "'8.17"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("E25").Value = '  -2.47%  '
$ws.Range("D26").Value = "'7.16"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = "'426.26"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("D31").Value = "'158.41"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  -1.57%  '
$ws.Range("D35").Value = "'17.74"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("D37").Value = "'4.39"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = "'129.47"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = "'0.0912"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  -3.21%  '
$ws.Range("D49").Value = "'4.92"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  -8.13%  '
$ws.Range("D50").Value = "'16.70"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  -4.69%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0204'
$ws.Range("E51").Value = '  -0.86%  '
